# Update "想去人数" (column F) values across the four sheets to reflect
# newly generated output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        4  = 1194
        13 = 115
        15 = 480
        35 = 241
    }
    "演出" = @{
        4  = 727
        15 = 380
        16 = 380
        19 = 933
    }
    "本地生活" = @{
        5  = 2285
        9  = 1139
        11 = 78
    }
    "全部类型" = @{
        4  = 2285
        10 = 1139
        12 = 78
        13 = 1194
        23 = 115
        24 = 480
        39 = 380
        51 = 241
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
